$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Cells whose new D value still looks numeric must keep their original
# "General" styling (no explicit style index), so we briefly switch the
# cell to Text format while writing the literal string, then restore the
# original style -- this avoids Excel auto-converting the text into a
# real number while also avoiding leaving a new number-format style behind.

$ws.Range("D2").Value = '37.385.38'
$ws.Range("E2").Value = '  +0.17%  '

$ws.Range("D3").Value = '2.070.55'
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("E4").Value = '  -0.04%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '235.02'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +0.36%  '

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.624'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  +2.14%  '

$ws.Range("E7").Value = '  -0.02%  '

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '57.20'
$cell.Style = $origStyle
$ws.Range("E8").Value = '  -1.43%  '

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.393'
$cell.Style = $origStyle
$ws.Range("E9").Value = '  +2.94%  '

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0775'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  +1.91%  '

$ws.Range("E11").Value = '  +0.92%  '

$ws.Range("D12").Value = '2.375.72'
$ws.Range("E12").Value = '  +0.44%  '

$ws.Range("E13").Value = '  -1.14%  '

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '20.81'
$cell.Style = $origStyle
$ws.Range("E14").Value = '  -1.50%  '

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.780'
$cell.Style = $origStyle
$ws.Range("E15").Value = '  +0.90%  '

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.20'
$cell.Style = $origStyle
$ws.Range("E16").Value = '  +0.45%  '

$ws.Range("D17").Value = '2.071.25'
$ws.Range("E17").Value = '  +0.23%  '

$ws.Range("D18").Value = '37.350.92'
$ws.Range("E18").Value = '  -0.51%  '

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.24'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  +1.28%  '

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '69.61'
$cell.Style = $origStyle
$ws.Range("E20").Value = '  +0.70%  '

$ws.Range("D21").Value = '0.0₃0819'
$ws.Range("E21").Value = '  +0.83%  '

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '226.68'
$cell.Style = $origStyle
$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("E24").Value = '  +1.58%  '

$ws.Range("E25").Value = '  -0.92%  '

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '168.00'
$cell.Style = $origStyle
$ws.Range("E26").Value = '  +1.66%  '

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.86'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("E28").Value = '  -5.65%  '

$ws.Range("E29").Value = '  +1.88%  '

$ws.Range("E30").Value = '  -0.24%  '

$ws.Range("E31").Value = '  -0.49%  '

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.55'
$cell.Style = $origStyle
$ws.Range("E32").Value = '  +1.23%  '

$ws.Range("E33").Value = '  -0.76%  '

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.53'
$cell.Style = $origStyle
$ws.Range("E34").Value = '  -0.55%  '

$ws.Range("E35").Value = '  -2.66%  '

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.37'
$cell.Style = $origStyle
$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("E37").Value = '  -1.00%  '

$ws.Range("E38").Value = '  -0.15%  '

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.63'
$cell.Style = $origStyle
$ws.Range("E39").Value = '  -3.72%  '

$ws.Range("E40").Value = '  -0.22%  '

$ws.Range("D41").Value = '1.491.53'
$ws.Range("E41").Value = '  +2.43%  '

$ws.Range("E42").Value = '  -3.13%  '

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '96.83'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  +1.45%  '

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0212'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  +1.11%  '

$ws.Range("E45").Value = '  +0.81%  '

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.11'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  -6.37%  '

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.04'
$cell.Style = $origStyle
$ws.Range("E47").Value = '  +0.93%  '

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '15.21'
$cell.Style = $origStyle
$ws.Range("E48").Value = '  -3.64%  '

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.22'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  -0.39%  '

$ws.Range("E50").Value = '  +0.71%  '

$ws.Range("D51").Value = '2.262.69'
$ws.Range("E51").Value = '  +0.45%  '
